$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.169.09"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.379.99"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.46"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.36"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +8.84%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.69"
$ws.Range("E11").Value = "  +4.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000286"
$ws.Range("E12").Value = "  +5.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "686.97"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.61"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.921.61"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.138.28"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.394.63"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.35"
$ws.Range("E20").Value = "  +3.18%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.07"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "104.01"
$ws.Range("E24").Value = "  +5.30%  "
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.66"
$ws.Range("E32").Value = "  +10.29%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.32"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.96"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.701.77"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").Value = "  +7.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.82"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0704"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.22"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0417"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.19"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  -2.14%  "
